$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text so numeric-looking strings
# (e.g. "3.500", "29.130.47") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.130.47'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').Value = '1.996.86'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('D4').Value = '1.014'
$ws.Range('E4').Value = '  +0.57%  '
$ws.Range('D5').Value = '330.93'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').Value = '0.4972'
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('D8').Value = '0.4193'
$ws.Range('E8').Value = '  -1.74%  '
$ws.Range('D9').Value = '54.93'
$ws.Range('E9').Value = '  +2.32%  '
$ws.Range('D10').Value = '0.08872'
$ws.Range('E10').Value = '  -3.68%  '
$ws.Range('D11').Value = '1.095'
$ws.Range('E11').Value = '  -3.15%  '
$ws.Range('D12').Value = '22.97'
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('D13').Value = '1.999.88'
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('D14').Value = '7.989'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').Value = '1.014'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '92.57'
$ws.Range('E17').Value = '  -3.61%  '
$ws.Range('D18').Value = '0.00001106'
$ws.Range('E18').Value = '  -1.86%  '
$ws.Range('D19').Value = '0.06754'
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').Value = '19.49'
$ws.Range('E20').Value = '  -2.17%  '
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('D23').Value = '29.165.99'
$ws.Range('E23').Value = '  -1.66%  '
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').Value = '2.292'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').Value = '2.252.82'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').Value = '20.82'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').Value = '157.09'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('D29').Value = '6.278'
$ws.Range('E29').Value = '  -4.75%  '
$ws.Range('D30').Value = '2.251'
$ws.Range('E30').Value = '  -4.40%  '
$ws.Range('D31').Value = '127.12'
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('D32').Value = '1.043'
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').Value = '0.09869'
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('D34').Value = '1.532'
$ws.Range('E34').Value = '  -4.35%  '
$ws.Range('D35').Value = '5.827'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').Value = '3.729'
$ws.Range('D37').Value = '0.02416'
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '1.315'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '9.082'
$ws.Range('E39').Value = '  -6.14%  '
$ws.Range('D40').Value = '0.06382'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('D41').Value = '0.6480'
$ws.Range('E41').Value = '  -1.73%  '
$ws.Range('D42').Value = '11.57'
$ws.Range('E42').Value = '  -2.23%  '
$ws.Range('E43').Value = '  -4.79%  '
$ws.Range('D44').Value = '1.012'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = '0.6189'
$ws.Range('E45').Value = '  -2.92%  '
$ws.Range('D46').Value = '1.355'
$ws.Range('E46').Value = '  +4.50%  '
$ws.Range('D47').Value = '13.28'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').Value = '2.175'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('D49').Value = '0.00000000345'
$ws.Range('E49').Value = '  +7.18%  '
$ws.Range('D50').Value = '3.500'
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').Value = '2.175'
$ws.Range('E51').Value = '  +7.96%  '
